$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds a new data row (row 3) to the sheet. Its values are an
# exact duplicate of the existing data row (row 2) - same venue, date,
# result, teams, batsman and stats. Copy row 2 down to row 3 so the new
# row picks up identical values *and* identical cell formatting/type
# (the numeric-looking totals must stay stored as text, matching the
# sheet's existing numberStoredAsText convention).
$ws.Range("A2:K2").Copy()
$ws.Range("A3").PasteSpecial()
